# Apply "Benar"/"Benaar" markers to column D and K for each of the 10
# test-sample blocks, then move the visible viewport / selection to
# match where the author ended up after testing all samples.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D markers ("apakah benar" flags) - one per sample block.
$ws.Range("D8").Value  = "Benar"
$ws.Range("D14").Value = "Benar"
$ws.Range("D22").Value = "Benar"
$ws.Range("D30").Value = "Benar"
$ws.Range("D38").Value = "Benar"
$ws.Range("D46").Value = "Benar"
$ws.Range("D54").Value = "Benar"
$ws.Range("D62").Value = "Benar"
$ws.Range("D69").Value = "Benaar"
$ws.Range("D77").Value = "Benar"

# Column K markers (mirrors column D, always "Benar" - bug fix row 69).
$ws.Range("K7").Value  = "Benar"
$ws.Range("K14").Value = "Benar"
$ws.Range("K22").Value = "Benar"
$ws.Range("K30").Value = "Benar"
$ws.Range("K38").Value = "Benar"
$ws.Range("K46").Value = "Benar"
$ws.Range("K54").Value = "Benar"
$ws.Range("K62").Value = "Benar"
$ws.Range("K69").Value = "Benar"
$ws.Range("K77").Value = "Benar"

# Leave the view scrolled down to the last edited sample and select it,
# matching where testing finished.
$ws.Range("D77").Select()
$excel.ActiveWindow.ScrollRow = 65
